$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 416) holds a "Förändrad" date that was bumped
# from 45206 (2023-10-07) to 45208 (2023-10-09) for every data row.
$ws.Range("C2:C416").Value = 45208
